# "yank debug prints from tools_xl.py"
#
# This edit:
#   1. Renames the provenance-table entry that records which python
#      script generated the workbook: "echo.pyc" -> "tools_xl.py".
#   2. Bumps the recorded generation "timestamp" to reflect the re-run.
#   3. Re-labels the three requirements summary cells ("of" -> "or").
#   4. Adds a new "08-BC" (Boundary Conditions) sheet at the end.
#   5. Reorders the tabs so "provenance" leads, followed by the three
#      requirements sheets, then the new "08-BC" sheet - and leaves
#      "provenance" as the selected/active tab.

$wb = $excel.ActiveWorkbook

# --- 1 & 2: update the provenance metadata table (do this before any
#     sheet reordering, while "provenance" is still easy to grab) ---
$prov = $wb.Worksheets.Item("provenance")
$prov.Range("B4").Value2 = "tools_xl.py"
$prov.Range("B12").Value2 = 43434.8592612651

# --- 3: fix the wording on the three requirements summary sheets ---
$wb.Worksheets.Item("requirements - PASS").Range("A1").Value2 = "Summary or requirements PASSED"
$wb.Worksheets.Item("requirements - FAIL").Range("A1").Value2 = "Summary or requirements FAIL"
$wb.Worksheets.Item("requirements - NULL").Range("A1").Value2 = "Summary or requirements NULL"

# --- 4: add the new "08-BC" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$bcSheet = $wb.Worksheets.Add($null, $lastSheet)
$bcSheet.Name = "08-BC"
$bcSheet.Range("A1").Value2 = "08-Boundary Conditions"

# --- 5: move "provenance" to be the first tab and make it the active
#     sheet. Re-fetch the reference by name (rather than reusing $prov)
#     since worksheet handles can go stale across a Move). ---
$wb.Worksheets.Item("provenance").Move($wb.Worksheets.Item(1))
$wb.Worksheets.Item("provenance").Activate()
